# The edit rotates the data in rows 2-4 (species/taxon records) upward:
#   new row 2 = old row 3
#   new row 3 = old row 4
#   new row 4 = old row 2
# Only the columns below actually differ between the three rows; all other
# columns (C, I, S, T, U, V, W, Y, Z, AA, AB, AD, AE, AG, AT, AW, AX, AY, ...)
# are identical across rows 2-4 and therefore remain untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ----- old values (as found in the workbook before this edit) -----
$row2 = @{
    A  = 74083298
    B  = 101808
    D  = "NT"
    E  = 223825
    F  = "Åkerkål"
    G  = "Brassica rapa subsp. campestris"
    H  = "(L.) A. R. Clapham"
    P  = "19 NO Strömserum, Sm"
    Q  = 561551.4394997598
    R  = 6326988.992090961
    AC = "Smålands flora 2007: KOO: 5G5c 3427. SOM: Brassica rapa ssp. campestris. LEG: Jörgen Andersson"
    AI = "Åker"
}

$row3 = @{
    A  = 74302538
    B  = 96309
    D  = "LC"
    E  = 219797
    F  = "Purpurknipprot"
    G  = "Epipactis atrorubens"
    H  = "(Hoffm.) Besser"
    P  = "18 S Sinnerboviken, Sm"
    Q  = 561645.2683258571
    R  = 6327489.84356571
    AC = "Smålands flora 2007: KOO: 5G5c 3928. SOM: Epipactis atrorubens. LEG: Jörgen Andersson"
    AI = "Granskog"
}

$row4 = @{
    A  = 74948080
    B  = 108194
    D  = "LC"
    E  = 219711
    F  = "Sårläka"
    G  = "Sanicula europaea"
    H  = "L."
    P  = "18 S Sinnerboviken, Sm"
    Q  = 561645.2683258571
    R  = 6327489.84356571
    AC = "Smålands flora 2007: KOO: 5G5c 3928. SOM: Sanicula europaea. LEG: Jörgen Andersson"
    AI = "Översilad skogsbacke"
}

# ----- apply rotation: row2 <- row3, row3 <- row4, row4 <- row2 -----
function Set-RowValues($targetRow, $values) {
    foreach ($col in $values.Keys) {
        $ws.Range("$col$targetRow").Value = $values[$col]
    }
}

Set-RowValues 2 $row3
Set-RowValues 3 $row4
Set-RowValues 4 $row2
